$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(187, 'Quaker Oats (1000gm)', 477386831.61, 29716069.6, 13645.134, 59690),
    @(172, 'Lays (52 gm) - STT', 1363305138, 20027118.3, 848663.634, 25148.76),
    @(170, 'Lays(52 gm) ASCO', 798975810, 15161443.3, 953873.036, 19038.76),
    @(195, 'Kurkure Large (45 gm) - ASCO', 140045089.2, 2089866.24, 153289.872, 4740.12),
    @(197, 'Kurkure Large (45 gm) - MM', 55874725.28, 1704890.88, 62366.172, 3866.94),
    @(168, 'Lays (25 gm) - STT', 79585686, 1615014, 101235.432, 1950.5),
    @(166, 'Lays(25 gm) ASCO', 126078225, 1615014, 137464.6, 1950.5),
    @(196, 'Kurkure Large (45 gm) - CC', 80504617.2, 1264919.04, 35556.528, 2869.02),
    @(213, 'Lays 3D-37 gm (Poly)', 18093323.04, 1204988.4, 78005.928, 2621.08),
    @(216, 'Lays Pastazz - 20 gm', 113853749.4, 942310, 341405.42, 1995.48),
    @(159, 'Kurkure (45 gm) - STT', 27437588.45, 895180.8, 50581.056, 2030.4),
    @(146, 'Kurkure (20 gm)- ASCO', 353070781.84, 828459.6, 672112.1040000001, 2086.8),
    @(142, 'Kurkure (20 gm) - Chilli Chatka', 144274267.2, 716505.6, 221693.952, 1804.8),
    @(160, 'Kurkure (90 gm) - STT', 1671383.7, 671385.6, 2887.92, 1522.8),
    @(215, 'Lays Pastazz - 37 gm (Poly)', 61251584.4, 641617.2, 62493.84, 1395.64),
    @(143, 'Kurkure (20 gm) - Masala Munch', 234505493.6, 559770, 222252, 1410),
    @(199, 'Quaker Oats (500 gm Jar)', 59237349.15, 549001.53, 384.561, 1039.5),
    @(158, 'Kurkure (20 gm) - STT', 148290535.05, 425425.2, 270590.172, 1071.6),
    @(155, 'Kurkure (90 gm) - Masala Munch', 2020411.62, 335692.8, 872.748, 761.4),
    @(157, 'Kurkure (90 gm) - ASCO', 6362383.95, 307718.4, 1679.535, 697.95),
    @(154, 'Kurkure (90 gm) - Chilli Chatka', 1620212.94, 251769.6, 524.907, 571.05),
    @(198, 'Kurkure Large (45 gm) - TJM', 1935023.68, 219985.92, 1114.752, 498.96),
    @(144, 'Kurkure (20 gm) - Tock Misti Jhal', 5189623.32, 134344.8, 7476.048, 338.4),
    @(217, 'Lays 3D - 20 gm', 5719548.91, 13005, 6404.87, 27.54),
    @('', '', 0, 0, 0, 0),
    @(156, 'Kurkure (90 gm) - Tock Misti Jhal', 593547.48, 0, 0, 0)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $data[$i]
    $r = $i + 2
    for ($j = 0; $j -lt $row.Count; $j++) {
        $val = $row[$j]
        if ($val -ne '') {
            $ws.Cells.Item($r, $j + 1).Value = $val
        }
    }
}
